$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 20
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 10000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 10000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -9770
$ws.Range("N20").ClearContents()

# Row 35
$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 10000
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 10000
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -9621
$ws.Range("N35").ClearContents()

# Row 58
$ws.Range("H58").Value = 50000050
$ws.Range("I58").Value = 50000050
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 150000150
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -150000000
$ws.Range("N58").ClearContents()

# Row 87
$ws.Range("I87").Value = 35000
$ws.Range("K87").Value = 35000
$ws.Range("M87").Value = -33752

# Row 90
$ws.Range("I90").Value = 35000
$ws.Range("K90").Value = 105000
$ws.Range("M90").Value = -98760

# Row 129
$ws.Range("H129").Value = 973.0833
$ws.Range("J129").Value = 1117.5
$ws.Range("L129").Value = 3352.5
$ws.Range("N129").Value = -13352.5

# Row 138
$ws.Range("H138").Value = 2617.74
$ws.Range("I138").Value = 1327.7354
$ws.Range("J138").Value = 3282.2878
$ws.Range("K138").Value = 3983.2062
$ws.Range("L138").Value = 9846.8634
$ws.Range("M138").Value = 1156.7938
$ws.Range("N138").Value = -20126.8634

# Row 139
$ws.Range("H139").Value = 248000
$ws.Range("J139").Value = 248000
$ws.Range("L139").Value = 248000
$ws.Range("N139").Value = -258280

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 21440.385
$ws.Range("I32").Value = 18667.162
$ws.Range("J32").Value = 72745
$ws.Range("K32").Value = 18667.162
$ws.Range("L32").Value = 72745
$ws.Range("M32").Value = -18380.162
$ws.Range("N32").Value = -73319

# Row 45
$ws.Range("H45").Value = 2088.4443
$ws.Range("I45").Value = 1514.7693
$ws.Range("J45").Value = 3580
$ws.Range("K45").Value = 1514.7693
$ws.Range("L45").Value = 3580
$ws.Range("M45").Value = -1137.7693
$ws.Range("N45").Value = -4334

# Row 63
$ws.Range("H63").Value = 2158.5
$ws.Range("I63").Value = 2176.111
$ws.Range("K63").Value = 2176.111
$ws.Range("M63").Value = -1490.111

# Row 66
$ws.Range("H66").Value = 2158.5
$ws.Range("I66").Value = 2176.111
$ws.Range("K66").Value = 10880.555
$ws.Range("M66").Value = -7448.555

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 8977
$ws.Range("I26").Value = 8977
$ws.Range("K26").Value = 8977
$ws.Range("M26").Value = -8685

# Row 86
$ws.Range("H86").Value = 2523.4092
$ws.Range("I86").Value = 2442.6667
$ws.Range("J86").Value = 2620.3
$ws.Range("K86").Value = 2442.6667
$ws.Range("L86").Value = 2620.3
$ws.Range("M86").Value = -1319.6667
$ws.Range("N86").Value = -4866.3

# Row 89
$ws.Range("H89").Value = 2523.4092
$ws.Range("I89").Value = 2442.6667
$ws.Range("J89").Value = 2620.3
$ws.Range("K89").Value = 12213.3335
$ws.Range("L89").Value = 13101.5
$ws.Range("M89").Value = -6597.333500000001
$ws.Range("N89").Value = -24333.5

# Row 96
$ws.Range("H96").Value = 9579.5
$ws.Range("I96").Value = 2949
$ws.Range("K96").Value = 2949
$ws.Range("M96").Value = -203

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4095.389
$ws.Range("I58").Value = 946
$ws.Range("J58").Value = 6913.263
$ws.Range("K58").Value = 946
$ws.Range("L58").Value = 6913.263
$ws.Range("M58").Value = -743
$ws.Range("N58").Value = -7319.263

# Row 115
$ws.Range("H115").Value = 27428.572
$ws.Range("J115").Value = 27428.572
$ws.Range("L115").Value = 27428.572
$ws.Range("N115").Value = -29778.572

# Row 136
$ws.Range("H136").Value = 4095.389
$ws.Range("I136").Value = 946
$ws.Range("J136").Value = 6913.263
$ws.Range("K136").Value = 2838
$ws.Range("L136").Value = 20739.789
$ws.Range("M136").Value = -288
$ws.Range("N136").Value = -25839.789

# Row 140
$ws.Range("H140").Value = 46640
$ws.Range("J140").Value = 46640
$ws.Range("L140").Value = 46640
$ws.Range("N140").Value = -57000

$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 2132.6667
$ws.Range("I34").Value = 900
$ws.Range("J34").Value = 3365.3333
$ws.Range("K34").Value = 2700
$ws.Range("L34").Value = 10095.9999
$ws.Range("M34").Value = -2616
$ws.Range("N34").Value = -10263.9999

# Row 39
$ws.Range("H39").Value = 5150
$ws.Range("J39").Value = 5150
$ws.Range("L39").Value = 15450
$ws.Range("N39").Value = -16038

# Row 55
$ws.Range("H55").Value = 40997.64
$ws.Range("I55").Value = 500400
$ws.Range("J55").Value = 1049.6086
$ws.Range("K55").Value = 1501200
$ws.Range("L55").Value = 3148.8258
$ws.Range("M55").Value = -1501023
$ws.Range("N55").Value = -3502.8258

# Row 68
$ws.Range("H68").Value = 978.7
$ws.Range("I68").Value = 508.14035
$ws.Range("J68").Value = 1602.4651
$ws.Range("K68").Value = 1524.42105
$ws.Range("L68").Value = 4807.3953
$ws.Range("M68").Value = -713.4210499999999
$ws.Range("N68").Value = -6429.3953

# Row 70
$ws.Range("H70").Value = 3874.75
$ws.Range("I70").Value = 4166.3335
$ws.Range("K70").Value = 12499.0005
$ws.Range("M70").Value = -12184.0005

# Row 71
$ws.Range("H71").Value = 978.7
$ws.Range("I71").Value = 508.14035
$ws.Range("J71").Value = 1602.4651
$ws.Range("K71").Value = 4573.26315
$ws.Range("L71").Value = 14422.1859
$ws.Range("M71").Value = -517.2631499999998
$ws.Range("N71").Value = -22534.1859

# Row 73
$ws.Range("H73").Value = 3874.75
$ws.Range("I73").Value = 4166.3335
$ws.Range("K73").Value = 12499.0005
$ws.Range("M73").Value = -11407.0005

# Row 131
$ws.Range("H131").Value = 24148.91
$ws.Range("I131").Value = 111633.336
$ws.Range("J131").Value = 2277.8057
$ws.Range("K131").Value = 334900.008
$ws.Range("L131").Value = 6833.4171
$ws.Range("M131").Value = -329860.008
$ws.Range("N131").Value = -16913.4171

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 1445.3077
$ws.Range("I113").Value = 1169.8889
$ws.Range("K113").Value = 1169.8889
$ws.Range("M113").Value = 1000.1111

# Row 138
$ws.Range("H138").Value = 39814.285
$ws.Range("J138").Value = 39814.285
$ws.Range("L138").Value = 39814.285
$ws.Range("N138").Value = -50094.285

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 2100.5
$ws.Range("I68").Value = 3000
$ws.Range("J68").Value = 1650.75
$ws.Range("K68").Value = 3000
$ws.Range("L68").Value = 1650.75
$ws.Range("M68").Value = -2251
$ws.Range("N68").Value = -3148.75

# Row 71
$ws.Range("H71").Value = 2100.5
$ws.Range("I71").Value = 3000
$ws.Range("J71").Value = 1650.75
$ws.Range("K71").Value = 15000
$ws.Range("L71").Value = 8253.75
$ws.Range("M71").Value = -11256
$ws.Range("N71").Value = -15741.75

# Row 138
$ws.Range("H138").Value = 35714.5
$ws.Range("J138").Value = 35714.5
$ws.Range("L138").Value = 35714.5
$ws.Range("N138").Value = -45994.5

# Row 139
$ws.Range("H139").Value = 48705
$ws.Range("J139").Value = 48705
$ws.Range("L139").Value = 48705
$ws.Range("N139").Value = -58985

$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 19896.75
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 19896.75
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 19896.75
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -20348.75

# Row 55
$ws.Range("H55").Value = 9487.25
$ws.Range("I55").Value = 1450
$ws.Range("K55").Value = 1450
$ws.Range("M55").Value = -1173

# Row 113
$ws.Range("H113").Value = 549
$ws.Range("I113").Value = 585.7143
$ws.Range("J113").Value = 463.33334
$ws.Range("K113").Value = 1757.1429
$ws.Range("L113").Value = 1390.00002
$ws.Range("M113").Value = 412.8571000000002
$ws.Range("N113").Value = -5730.000019999999

# Row 123
$ws.Range("H123").Value = 34529.69
$ws.Range("J123").Value = 34529.69
$ws.Range("L123").Value = 34529.69
$ws.Range("N123").Value = -44329.69

# Row 132
$ws.Range("H132").Value = 2584.1177
$ws.Range("I132").Value = 1134
$ws.Range("J132").Value = 3375.0908
$ws.Range("K132").Value = 3402
$ws.Range("L132").Value = 10125.2724
$ws.Range("M132").Value = -872
$ws.Range("N132").Value = -15185.2724

# Row 138
$ws.Range("H138").Value = 41896.668
$ws.Range("I138").Value = 30390
$ws.Range("J138").Value = 47650
$ws.Range("K138").Value = 30390
$ws.Range("L138").Value = 47650
$ws.Range("M138").Value = -25250
$ws.Range("N138").Value = -57930
